$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 6 (categoryMotility/Motility),
# shifting the existing rows 6-12 down to 8-14.
$ws.Rows("6:7").Insert()

# Populate the two newly-inserted rows with the new category key/value pairs.
$ws.Range("A6").Value = "categoryMetabolism"
$ws.Range("B6").Value = "Metabolism"

$ws.Range("A7").Value = "categoryCoating"
$ws.Range("B7").Value = "Coating"

# Update the active selection to match the author's final cursor position.
$ws.Range("A15").Select()
